$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8: Assignment of real var to int var -> add Implemented=Y, Notes
$ws.Range("D8").Value = "Y"
$ws.Range("E8").Value = "Warning generated warns loss of precision"

# Row 9: Assignment of int var to char var -> add Implemented=N, Notes
$ws.Range("D9").Value = "N"
$ws.Range("E9").Value = "Warning generated"

# Row 10: Assignment of real var to char var -> add Implemented=Y
$ws.Range("D10").Value = "Y"

# Row 15: Write uninitialised var -> add Implemented=Y, Notes
$ws.Range("D15").Value = "Y"
$ws.Range("E15").Value = "All variables are initilised to one."

# Row 17: Extra semicolon inside if -> add Implemented=Y, Notes
$ws.Range("D17").Value = "Y"
$ws.Range("E17").Value = "Warns and ignores"

# Row 20: fix typo (remove leading space) in description
$ws.Range("C20").Value = "does compiler prevent constant divide by zero"

# Row 21: invalid char input -> add Notes
$ws.Range("E21").Value = "Yes removes all until read a valid string. Removes white space. Flushes after reading."

# New rows 23-25: additional "Me" notes
$ws.Range("A23").Value = "?"
$ws.Range("B23").Value = "Me"
$ws.Range("C23").Value = "Constant folding"

$ws.Range("A24").Value = "?"
$ws.Range("B24").Value = "Me"
$ws.Range("C24").Value = "Loop unwinding"

$ws.Range("A25").Value = "?"
$ws.Range("B25").Value = "Me"
$ws.Range("C25").Value = "Redundant assignment removal"

# Update column E width to fit new, longer notes (bestFit-style autosize)
$ws.Columns.Item(5).ColumnWidth = 70

# Update selection to match final state (active cell C25)
$ws.Range("C25").Select()
